$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 401.1111
$ws.Range("H113").Value = 2646.7942
$ws.Range("I113").Value = 2712.5
$ws.Range("J113").Value = 2572.875
$ws.Range("K113").Value = 2712.5
$ws.Range("L113").Value = 2572.875
$ws.Range("M113").Value = 541.5
$ws.Range("N113").Value = -9080.875
$ws.Range("H137").Value = 11112210
$ws.Range("I137").Value = 14706753
$ws.Range("J137").Value = 1805.7273
$ws.Range("K137").Value = 44120259
$ws.Range("L137").Value = 5417.1819
$ws.Range("M137").Value = -44117709
$ws.Range("N137").Value = -10517.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 1063
$ws.Range("I26").Value = 826.75
$ws.Range("J26").Value = 2008
$ws.Range("K26").Value = 826.75
$ws.Range("L26").Value = 2008
$ws.Range("M26").Value = -496.75
$ws.Range("N26").Value = -2668
$ws.Range("H61").Value = 3334695
$ws.Range("I61").Value = 4066323
$ws.Range("J61").Value = 1722.2222
$ws.Range("K61").Value = 4066323
$ws.Range("L61").Value = 1722.2222
$ws.Range("M61").Value = -4066111
$ws.Range("N61").Value = -2146.2222
$ws.Range("H74").Value = 1103.919
$ws.Range("I74").Value = 1025.6061
$ws.Range("K74").Value = 1025.6061
$ws.Range("M74").Value = -151.6061
$ws.Range("H77").Value = 1103.919
$ws.Range("I77").Value = 1025.6061
$ws.Range("K77").Value = 5128.0305
$ws.Range("M77").Value = -760.0304999999998
$ws.Range("H132").Value = 806921.9
$ws.Range("I132").Value = 897.06665
$ws.Range("K132").Value = 2691.19995
$ws.Range("M132").Value = -161.1999500000002
$ws.Range("H136").Value = 3334695
$ws.Range("I136").Value = 4066323
$ws.Range("J136").Value = 1722.2222
$ws.Range("K136").Value = 12198969
$ws.Range("L136").Value = 5166.6666
$ws.Range("M136").Value = -12196419
$ws.Range("N136").Value = -10266.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1293766.1
$ws.Range("I86").Value = 1497.6666
$ws.Range("J86").Value = 3878303.2
$ws.Range("K86").Value = 1497.6666
$ws.Range("L86").Value = 3878303.2
$ws.Range("M86").Value = -374.6666
$ws.Range("N86").Value = -3880549.2
$ws.Range("H89").Value = 1293766.1
$ws.Range("I89").Value = 1497.6666
$ws.Range("J89").Value = 3878303.2
$ws.Range("K89").Value = 7488.333000000001
$ws.Range("L89").Value = 19391516
$ws.Range("M89").Value = -1872.333000000001
$ws.Range("N89").Value = -19402748
$ws.Range("H107").Value = 23810280
$ws.Range("I107").Value = 29412472
$ws.Range("J107").Value = 962
$ws.Range("K107").Value = 29412472
$ws.Range("L107").Value = 962
$ws.Range("M107").Value = -29410552
$ws.Range("N107").Value = -4802
$ws.Range("H134").Value = 2181827
$ws.Range("I134").Value = 1181.6342
$ws.Range("J134").Value = 11122473
$ws.Range("K134").Value = 3544.9026
$ws.Range("L134").Value = 33367419
$ws.Range("M134").Value = -1009.9026
$ws.Range("N134").Value = -33372489

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1303.4849
$ws.Range("I16").Value = 801.06665
$ws.Range("J16").Value = 1722.1666
$ws.Range("K16").Value = 801.06665
$ws.Range("L16").Value = 1722.1666
$ws.Range("M16").Value = -514.06665
$ws.Range("N16").Value = -2296.1666
$ws.Range("H31").Value = 928620.25
$ws.Range("I31").Value = 1184959.9
$ws.Range("J31").Value = 1853.8462
$ws.Range("K31").Value = 1184959.9
$ws.Range("L31").Value = 1853.8462
$ws.Range("M31").Value = -1184664.9
$ws.Range("N31").Value = -2443.8462
$ws.Range("H34").Value = 928620.25
$ws.Range("I34").Value = 1184959.9
$ws.Range("J34").Value = 1853.8462
$ws.Range("K34").Value = 1184959.9
$ws.Range("L34").Value = 1853.8462
$ws.Range("M34").Value = -1184757.9
$ws.Range("N34").Value = -2257.8462
$ws.Range("H58").Value = 32258850
$ws.Range("I58").Value = 71429590
$ws.Range("J58").Value = 594.35297
$ws.Range("K58").Value = 71429590
$ws.Range("L58").Value = 594.35297
$ws.Range("M58").Value = -71429387
$ws.Range("N58").Value = -1000.35297
$ws.Range("H113").Value = 1303.4849
$ws.Range("I113").Value = 801.06665
$ws.Range("J113").Value = 1722.1666
$ws.Range("K113").Value = 801.06665
$ws.Range("L113").Value = 1722.1666
$ws.Range("M113").Value = 1368.93335
$ws.Range("N113").Value = -6062.1666
$ws.Range("H132").Value = 6667933
$ws.Range("I132").Value = 927.26666
$ws.Range("J132").Value = 16668442
$ws.Range("K132").Value = 2781.79998
$ws.Range("L132").Value = 50005326
$ws.Range("M132").Value = -251.7999799999998
$ws.Range("N132").Value = -50010386
$ws.Range("H134").Value = 26317006
$ws.Range("I134").Value = 1194.6154
$ws.Range("J134").Value = 83334600
$ws.Range("K134").Value = 3583.8462
$ws.Range("L134").Value = 250003800
$ws.Range("M134").Value = -1048.8462
$ws.Range("N134").Value = -250008870
$ws.Range("H136").Value = 32258850
$ws.Range("I136").Value = 71429590
$ws.Range("J136").Value = 594.35297
$ws.Range("K136").Value = 214288770
$ws.Range("L136").Value = 1783.05891
$ws.Range("M136").Value = -214286220
$ws.Range("N136").Value = -6883.05891

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 154.73334
$ws.Range("I6").Value = 101.90909
$ws.Range("K6").Value = 305.72727
$ws.Range("M6").Value = -192.72727
$ws.Range("H7").Value = 307.2857
$ws.Range("I7").Value = 307.2857
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 921.8571000000001
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -809.8571000000001
$ws.Range("N7").ClearContents()
$ws.Range("H92").Value = 7975.3335
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 7975.3335
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 23926.0005
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -26422.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9092854
$ws.Range("I80").Value = 2124.375
$ws.Range("J80").Value = 33334800
$ws.Range("K80").Value = 2124.375
$ws.Range("L80").Value = 33334800
$ws.Range("M80").Value = -1126.375
$ws.Range("N80").Value = -33336796
$ws.Range("H83").Value = 9092854
$ws.Range("I83").Value = 2124.375
$ws.Range("J83").Value = 33334800
$ws.Range("K83").Value = 10621.875
$ws.Range("L83").Value = 166674000
$ws.Range("M83").Value = -5629.875
$ws.Range("N83").Value = -166683984
$ws.Range("H132").Value = 3299.5757
$ws.Range("I132").Value = 1545.8334
$ws.Range("J132").Value = 11191.417
$ws.Range("K132").Value = 4637.5002
$ws.Range("L132").Value = 33574.251
$ws.Range("M132").Value = -2107.5002
$ws.Range("N132").Value = -38634.251

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 18187380
$ws.Range("I132").Value = 27028520
$ws.Range("J132").Value = 13926.833
$ws.Range("K132").Value = 81085560
$ws.Range("L132").Value = 41780.499
$ws.Range("M132").Value = -81083030
$ws.Range("N132").Value = -46840.499
$ws.Range("H136").Value = 39898644
$ws.Range("I136").Value = 19843124
$ws.Range("J136").Value = 76924220
$ws.Range("K136").Value = 59529372
$ws.Range("L136").Value = 230772660
$ws.Range("M136").Value = -59526822
$ws.Range("N136").Value = -230777760

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7594029
$ws.Range("I132").Value = 24994.674
$ws.Range("J132").Value = 21744832
$ws.Range("K132").Value = 74984.022
$ws.Range("L132").Value = 65234496
$ws.Range("M132").Value = -72454.022
$ws.Range("N132").Value = -65239556
$ws.Range("H136").Value = 17859632
$ws.Range("I136").Value = 38463350
$ws.Range("J136").Value = 3074.6667
$ws.Range("K136").Value = 115390050
$ws.Range("L136").Value = 9224.000100000001
$ws.Range("M136").Value = -115387500
$ws.Range("N136").Value = -14324.0001
